# Fruta / hortaliza, semanal
# Insert two new weekly price rows (Murcott Primera / Segunda, Terminal
# Hortofrutícola Agro Chillán, Provincia de Limarí, 2021-10-13) right before
# the existing row 55, shifting the rest of the table down by two rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push existing rows 55.. down by two rows, creating two blank rows at 55:56
# (Excel copies the formatting of the row above the insertion point, which
# keeps the date-formatted column D intact for the new rows too).
$ws.Rows("55:56").Insert()

$newDate = Get-Date -Year 2021 -Month 10 -Day 13 -Hour 0 -Minute 0 -Second 0

# Row 55: Murcott - Primera
$ws.Range("A55").Value = 7
$ws.Range("B55").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C55").Value = "Ñuble"
$ws.Range("D55").Value = $newDate
$ws.Range("E55").Value = 16
$ws.Range("F55").Value = "Fruta"
$ws.Range("G55").Value = 100102
$ws.Range("H55").Value = "Cítricos"
$ws.Range("I55").Value = 100102004
$ws.Range("J55").Value = "Mandarina"
$ws.Range("K55").Value = "Murcott"
$ws.Range("L55").Value = "Primera"
$ws.Range("M55").Value = 240
$ws.Range("N55").Value = 6000
$ws.Range("O55").Value = 6500
$ws.Range("P55").Value = 6250
$ws.Range("Q55").Value = "`$/bandeja 10 kilos"
$ws.Range("R55").Value = "Provincia de Limarí"
$ws.Range("S55").Value = 625
$ws.Range("T55").Value = 10

# Row 56: Murcott - Segunda
$ws.Range("A56").Value = 7
$ws.Range("B56").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C56").Value = "Ñuble"
$ws.Range("D56").Value = $newDate
$ws.Range("E56").Value = 16
$ws.Range("F56").Value = "Fruta"
$ws.Range("G56").Value = 100102
$ws.Range("H56").Value = "Cítricos"
$ws.Range("I56").Value = 100102004
$ws.Range("J56").Value = "Mandarina"
$ws.Range("K56").Value = "Murcott"
$ws.Range("L56").Value = "Segunda"
$ws.Range("M56").Value = 240
$ws.Range("N56").Value = 5000
$ws.Range("O56").Value = 5500
$ws.Range("P56").Value = 5250
$ws.Range("Q56").Value = "`$/bandeja 10 kilos"
$ws.Range("R56").Value = "Provincia de Limarí"
$ws.Range("S56").Value = 525
$ws.Range("T56").Value = 10
